$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.6
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 1.85
$ws.Range("U2").Value = 4.3
$ws.Range("V2").Value = 1.22
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 1.14
$ws.Range("AG2").Value = 21
$ws.Range("AO2").Value = 17
$ws.Range("G3").Value = 2.1
$ws.Range("I3").Value = 4.33
$ws.Range("K3").Value = 1.8
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("S3").Value = 3.4
$ws.Range("T3").Value = 1.33
$ws.Range("Y3").Value = 1.73
$ws.Range("Z3").Value = 2
$ws.Range("AE3").Value = 11
$ws.Range("H4").Value = 2.88
$ws.Range("L4").Value = 4
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 2.03
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 2.7
$ws.Range("T4").Value = 1.44
$ws.Range("U4").Value = 4.4
$ws.Range("V4").Value = 1.22
$ws.Range("W4").Value = 5.5
$ws.Range("X4").Value = 1.14
$ws.Range("Y4").Value = 1.62
$ws.Range("Z4").Value = 2.2
$ws.Range("AA4").Value = 2.1
$ws.Range("AB4").Value = 1.67
$ws.Range("AG4").Value = 26
$ws.Range("AI4").Value = 6
$ws.Range("AK4").Value = 19
$ws.Range("AN4").Value = 7
$ws.Range("G5").Value = 2.75
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 3.75
$ws.Range("AD5").Value = 11
$ws.Range("AF5").Value = 29
$ws.Range("AR5").Value = 41
$ws.Range("G6").Value = 1.91
$ws.Range("I6").Value = 5.25
$ws.Range("J6").Value = 2.75
$ws.Range("K6").Value = 1.83
$ws.Range("L6").Value = 6
$ws.Range("AC6").Value = 4.5
$ws.Range("AD6").Value = 7
$ws.Range("AF6").Value = 15
$ws.Range("AK6").Value = 26
$ws.Range("AL6").Value = 126
$ws.Range("AN6").Value = 9
$ws.Range("AO6").Value = 23
$ws.Range("AP6").Value = 21
$ws.Range("S7").Value = 3.4
$ws.Range("T7").Value = 1.33
$ws.Range("W7").Value = 7
$ws.Range("X7").Value = 1.1
$ws.Range("Z7").Value = 2.08
$ws.Range("G8").Value = 2.45
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 2.88
$ws.Range("J8").Value = 3.2
$ws.Range("L8").Value = 3.6
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("S8").Value = 2.08
$ws.Range("T8").Value = 1.73
$ws.Range("AC8").Value = 7.5
$ws.Range("AD8").Value = 11
$ws.Range("AF8").Value = 23
$ws.Range("AG8").Value = 21
$ws.Range("AN8").Value = 8.5
$ws.Range("AQ8").Value = 29
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("S9").Value = 2.25
$ws.Range("T9").Value = 1.62
$ws.Range("G10").Value = 2.15
$ws.Range("H10").Value = 3.05
$ws.Range("I10").Value = 3.3
$ws.Range("J10").Value = 2.77
$ws.Range("K10").Value = 2.02
$ws.Range("L10").Value = 3.95
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 5.9
$ws.Range("O10").Value = 1.45
$ws.Range("P10").Value = 2.57
$ws.Range("S10").Value = 2.3
$ws.Range("T10").Value = 1.55
$ws.Range("W10").Value = 4.05
$ws.Range("X10").Value = 1.2
$ws.Range("Y10").Value = 1.47
$ws.Range("Z10").Value = 2.5
$ws.Range("AA10").Value = 2
$ws.Range("AB10").Value = 1.72
$ws.Range("AC10").Value = 6.2
$ws.Range("AD10").Value = 9.25
$ws.Range("AE10").Value = 9.25
$ws.Range("AF10").Value = 20
$ws.Range("AG10").Value = 20
$ws.Range("AH10").Value = 37
$ws.Range("AI10").Value = 5.9
$ws.Range("AJ10").Value = 6.1
$ws.Range("AK10").Value = 17.5
$ws.Range("AL10").Value = 100
$ws.Range("AM10").Value = 900
$ws.Range("AN10").Value = 7.9
$ws.Range("AO10").Value = 16
$ws.Range("AP10").Value = 12.5
$ws.Range("AR10").Value = 35
$ws.Range("AS10").Value = 50
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 3.6
$ws.Range("S11").Value = 1.85
$ws.Range("T11").Value = 1.95
$ws.Range("AF11").Value = 17
$ws.Range("AG11").Value = 15
$ws.Range("AJ11").Value = 7
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.5
$ws.Range("Q12").Value = 1.92
$ws.Range("R12").Value = 1.82
$ws.Range("G13").Value = 2.6
$ws.Range("I13").Value = 2.6
$ws.Range("J13").Value = 3.4
$ws.Range("AC13").Value = 8.5
$ws.Range("AD13").Value = 13
$ws.Range("AG13").Value = 23
$ws.Range("AP13").Value = 10
